$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the stats for 2026-01 (row 26) to their new values
$ws.Range("B26").Value = 6532
$ws.Range("C26").Value = 1016
$ws.Range("D26").Value = 6093296
$ws.Range("E26").Value = 932.8377219840784
$ws.Range("F26").Value = 10.39378063207708
$ws.Range("G26").Value = 7.855626326963905
$ws.Range("H26").Value = 26.89309448130359
